$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the query text in column B for each row (order matters for
# shared-string table layout: participants, samples, then files)
$ws.Range("B2").Value = 'MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.reference_genome_assembly in [''GRCh37''] 
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'''') as `Participant ID`,
coalesce(s.study_name, '''') as `Study Name`,
coalesce(s.phs_accession,'''') as `Accession`,
coalesce(p.gender,'''') as `Gender`,
coalesce(apoc.text.join(samp, '',''), '''') as `Samples`
ORDER BY p.participant_id limit 100'
$ws.Range("B3").Value = 'MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.reference_genome_assembly in [''GRCh37'']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(p.participant_id,'''') as `Participant ID`,
    coalesce(s.study_name, '''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(samp.sample_tumor_status,'''') as `Tumor`,
    coalesce(samp.sample_type,'''') as `Analyte Type`
ORDER BY samp.sample_id limit 100'
$ws.Range("B4").Value = 'MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.reference_genome_assembly in [''GRCh37'']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '''') as `File Name`,
    coalesce(s.study_name,'''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(p.participant_id, '''') as `Participant ID`,
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(f.file_type, '''') as `File Type`
ORDER BY f.file_name limit 100'

# Update the stat query text in column C for each row
$ws.Range("C2").Value = 'CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in [''GRCh37'']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in [''GRCh37'']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in [''GRCh37'']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`'
$ws.Range("C3").Value = 'CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in [''GRCh37'']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in [''GRCh37'']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in [''GRCh37'']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`'
$ws.Range("C4").Value = 'CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in [''GRCh37'']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in [''GRCh37'']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in [''GRCh37'']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`'

# Re-pin the (wrap-text) row heights - writing the new, longer query text
# would otherwise trigger Excel's auto-fit and grow the rows
$ws.Rows.Item(2).RowHeight = 242.25
$ws.Rows.Item(3).RowHeight = 260.25
$ws.Rows.Item(4).RowHeight = 279.75

# Update the selection (also clears the scrolled-away topLeftCell)
$ws.Range("B2:C4").Select() | Out-Null
